$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell "D2" "62.306.07"
Set-TextCell "E2" "  +0.92%  "

Set-TextCell "D3" "3.428.11"
Set-TextCell "E3" "  +0.67%  "

Set-TextCell "E4" "  -0.27%  "

Set-TextCell "D5" "413.12"
Set-TextCell "E5" "  +0.89%  "

Set-TextCell "D6" "128.72"
Set-TextCell "E6" "  +0.12%  "

Set-TextCell "E7" "  -1.73%  "

Set-TextCell "E8" "  +0.00%  "

Set-TextCell "D9" "0.727"
Set-TextCell "E9" "  -0.77%  "

Set-TextCell "D10" "0.141"
Set-TextCell "E10" "  +0.93%  "

Set-TextCell "D11" "42.74"
Set-TextCell "E11" "  +0.32%  "

Set-TextCell "D12" "0.0000222"
Set-TextCell "E12" "  +7.40%  "

Set-TextCell "E13" "  +2.59%  "

Set-TextCell "D14" "3.971.51"
Set-TextCell "E14" "  +0.50%  "

Set-TextCell "E15" "  -0.13%  "

Set-TextCell "E16" "  -2.87%  "

Set-TextCell "D17" "3.433.55"
Set-TextCell "E17" "  +1.51%  "

Set-TextCell "D18" "12.71"
Set-TextCell "E18" "  +5.32%  "

Set-TextCell "E19" "  -0.24%  "

Set-TextCell "D20" "62.306.51"
Set-TextCell "E20" "  +0.89%  "

Set-TextCell "D21" "477.85"
Set-TextCell "E21" "  +7.37%  "

Set-TextCell "D22" "91.63"
Set-TextCell "E22" "  +0.54%  "

Set-TextCell "E23" "  +3.50%  "

Set-TextCell "E24" "  +1.55%  "

Set-TextCell "D25" "3.29"
Set-TextCell "E25" "  +1.28%  "

Set-TextCell "D26" "9.77"
Set-TextCell "E26" "  +12.11%  "

Set-TextCell "D27" "33.49"
Set-TextCell "E27" "  +0.37%  "

Set-TextCell "E28" "  +1.10%  "

Set-TextCell "D29" "7.77"
Set-TextCell "E29" "  +2.70%  "

Set-TextCell "D30" "11.87"
Set-TextCell "E30" "  -0.55%  "

Set-TextCell "E31" "  -3.27%  "

Set-TextCell "E32" "  -1.67%  "

Set-TextCell "E33" "  -1.61%  "

Set-TextCell "D34" "40.66"
Set-TextCell "E34" "  -4.37%  "

Set-TextCell "D36" "57.98"
Set-TextCell "E36" "  +8.05%  "

Set-TextCell "E37" "  -1.52%  "

Set-TextCell "D38" "1.00"
Set-TextCell "E38" "  +0.10%  "

Set-TextCell "D39" "3.03"
Set-TextCell "E39" "  +4.79%  "

Set-TextCell "E40" "  +1.05%  "

Set-TextCell "D41" "0.323"
Set-TextCell "E41" "  +3.51%  "

Set-TextCell "E42" "  -1.32%  "

Set-TextCell "B43" "WEMIXToken"
Set-TextCell "C43" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell "D43" "2.67"
Set-TextCell "E43" "  +11.50%  "

Set-TextCell "B44" "Monero"
Set-TextCell "C44" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D44" "144.65"
Set-TextCell "E44" "  +2.56%  "

Set-TextCell "D45" "4.32"
Set-TextCell "E45" "  +4.34%  "

Set-TextCell "E46" "  +4.49%  "

Set-TextCell "E47" "  +19.83%  "

Set-TextCell "E48" "  -1.52%  "

Set-TextCell "D49" "0.0₃0544"
Set-TextCell "E49" "  +29.89%  "

Set-TextCell "D50" "22.16"
Set-TextCell "E50" "  -0.32%  "

Set-TextCell "D51" "112.91"
Set-TextCell "E51" "  +7.41%  "
